$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price (D) column writes to remain text (matches source data which
# stores prices as plain text, even when they look numeric), then restore
# the cell style so no stray number-format style is left behind.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '88.016.43'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.25%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.106.16'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -2.60%  '

$ws.Range('E4').Value = '  +0.12%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.36'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.09%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '634.06'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.72%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.387'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.71%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.783'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +13.51%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '1.00'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.18%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.105.09'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.39%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.559'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.24%  '

$ws.Range('E12').Value = '  +0.91%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000248'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.83%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.36'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +2.18%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '88.216.39'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -2.49%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.693.79'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.80%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '31.98'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.72%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.120.56'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.30%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '3.37'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.13%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0000218'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +15.64%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.14'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.65%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '420.07'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.90%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.37'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -2.53%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '4.87'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.65%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '5.39'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +4.74%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '82.09'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +9.78%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.39'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -3.23%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '3.283.71'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.89%  '

$ws.Range('E29').Value = '  +0.04%  '

$ws.Range('E30').Value = '  +0.32%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.155'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -8.76%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.98'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.01%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '8.12'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.96%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.147'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +14.66%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '498.87'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.54%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.85'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.63%  '

$ws.Range('B37').Value = 'PancakeSwap'
$ws.Range('C37').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.83'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -2.08%  '

$ws.Range('B38').Value = 'Fetch.AI'
$ws.Range('C38').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.26'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.50%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '22.18'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.00%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '22.20'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.57%  '

$ws.Range('E41').Value = '  +0.39%  '

$ws.Range('E42').Value = '  +0.05%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.362'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.55%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.83'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.14%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '145.80'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.14%  '

$ws.Range('E46').Value = '  +6.92%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '43.48'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.12%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0652'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +10.83%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '161.88'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.23%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.712'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.80%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.18'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -4.84%  '
